# Auto-generated edit script applying the Golem_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets("ALC")
$ws.Range("H2").Value = 794.8333
$ws.Range("I2").Value = 843.5
$ws.Range("J2").Value = 770.5
$ws.Range("K2").Value = 843.5
$ws.Range("L2").Value = 770.5
$ws.Range("M2").Value = -730.5
$ws.Range("N2").Value = -996.5
$ws.Range("H13").Value = 1900
$ws.Range("J13").Value = 1900
$ws.Range("L13").Value = 1900
$ws.Range("N13").Value = -2238
$ws.Range("H29").Value = 425
$ws.Range("J29").Value = 800
$ws.Range("L29").Value = 2400
$ws.Range("N29").Value = -2962
$ws.Range("H38").Value = 1681.5
$ws.Range("I38").Value = 17.8
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 53.40000000000001
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = 318.6
$ws.Range("N38").Value = -30744
$ws.Range("H40").Value = 2900
$ws.Range("I40").Value = 2900
$ws.Range("K40").Value = 2900
$ws.Range("M40").Value = -2725
$ws.Range("H43").Value = 6000
$ws.Range("J43").Value = 8000
$ws.Range("L43").Value = 8000
$ws.Range("N43").Value = -8138
$ws.Range("H58").Value = 103
$ws.Range("J58").Value = 475
$ws.Range("L58").Value = 1425
$ws.Range("N58").Value = -1725
$ws.Range("H64").Value = 9999
$ws.Range("J64").Value = 9999
$ws.Range("L64").Value = 9999
$ws.Range("N64").Value = -10495
$ws.Range("H67").Value = 9999
$ws.Range("J67").Value = 9999
$ws.Range("L67").Value = 9999
$ws.Range("N67").Value = -11715
$ws.Range("H97").Value = 876.75
$ws.Range("J97").Value = 876.75
$ws.Range("L97").Value = 2630.25
$ws.Range("N97").Value = -3622.25
$ws.Range("H106").Value = 333333800
$ws.Range("I106").Value = 333333800
$ws.Range("K106").Value = 333333800
$ws.Range("M106").Value = -333333169
$ws.Range("H107").Value = 50637.223
$ws.Range("I107").Value = 53570.176
$ws.Range("K107").Value = 53570.176
$ws.Range("M107").Value = -51650.176
$ws.Range("H111").Value = 2339.7273
$ws.Range("I111").Value = 2339.7273
$ws.Range("K111").Value = 7019.1819
$ws.Range("M111").Value = -3952.1819
$ws.Range("H138").Value = 2814.625
$ws.Range("J138").Value = 3178.45
$ws.Range("L138").Value = 9535.349999999999
$ws.Range("N138").Value = -19815.35

# ---- Sheet: ARM ----
$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 1937.6666
$ws.Range("I32").Value = 1937.6666
$ws.Range("K32").Value = 1937.6666
$ws.Range("M32").Value = -1650.6666
$ws.Range("H61").Value = 3519.75
$ws.Range("I61").Value = 3519.75
$ws.Range("K61").Value = 3519.75
$ws.Range("M61").Value = -3307.75
$ws.Range("H101").Value = 27416.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 27416.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 27416.5
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -33906.5
$ws.Range("H132").Value = 1719.5454
$ws.Range("I132").Value = 1691.5
$ws.Range("K132").Value = 5074.5
$ws.Range("M132").Value = -2544.5
$ws.Range("H136").Value = 3519.75
$ws.Range("I136").Value = 3519.75
$ws.Range("K136").Value = 10559.25
$ws.Range("M136").Value = -8009.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets("BSM")
$ws.Range("H9").Value = 94306.336
$ws.Range("J9").Value = 94306.336
$ws.Range("L9").Value = 94306.336
$ws.Range("N9").Value = -94642.336
$ws.Range("H105").Value = 1925.6666
$ws.Range("I105").Value = 1925.6666
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1925.6666
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -178.6666
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 25682.938
$ws.Range("I107").Value = 29224.215
$ws.Range("K107").Value = 29224.215
$ws.Range("M107").Value = -27304.215

# ---- Sheet: CRP ----
$ws = $wb.Worksheets("CRP")
$ws.Range("H44").Value = 35999.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 35999.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 35999.5
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -36883.5
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H99").Value = 1429999.8
$ws.Range("I99").Value = 1251249.8
$ws.Range("J99").Value = 1668333
$ws.Range("K99").Value = 1251249.8
$ws.Range("L99").Value = 1668333
$ws.Range("M99").Value = -1249751.8
$ws.Range("N99").Value = -1671329
$ws.Range("H126").Value = 1429999.8
$ws.Range("I126").Value = 1251249.8
$ws.Range("J126").Value = 1668333
$ws.Range("K126").Value = 3753749.4
$ws.Range("L126").Value = 5004999
$ws.Range("M126").Value = -3751279.4
$ws.Range("N126").Value = -5009939
$ws.Range("H134").Value = 1308.3334
$ws.Range("I134").Value = 1221.875
$ws.Range("K134").Value = 3665.625
$ws.Range("M134").Value = -1130.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets("CUL")
$ws.Range("H114").Value = 2500
$ws.Range("I114").Value = 2500
$ws.Range("K114").Value = 7500
$ws.Range("M114").Value = -4246
$ws.Range("H132").Value = 659
$ws.Range("I132").Value = 559.8
$ws.Range("K132").Value = 5038.2
$ws.Range("M132").Value = -2508.2
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets("GSM")
$ws.Range("H58").Value = 30020.5
$ws.Range("I58").Value = 30020.5
$ws.Range("K58").Value = 30020.5
$ws.Range("M58").Value = -29743.5
$ws.Range("H102").Value = 782.5714
$ws.Range("I102").Value = 782.5714
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 782.5714
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 839.4286
$ws.Range("N102").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets("LTW")
$ws.Range("H46").Value = 1419.2858
$ws.Range("J46").Value = 1883.75
$ws.Range("L46").Value = 1883.75
$ws.Range("N46").Value = -2259.75
$ws.Range("H55").Value = 881.5
$ws.Range("I55").Value = 616.4
$ws.Range("J55").Value = 1544.25
$ws.Range("K55").Value = 616.4
$ws.Range("L55").Value = 1544.25
$ws.Range("M55").Value = -443.4
$ws.Range("N55").Value = -1890.25
$ws.Range("H136").Value = 4937
$ws.Range("I136").Value = 3998.6667
$ws.Range("K136").Value = 11996.0001
$ws.Range("M136").Value = -9446.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets("WVR")
$ws.Range("H52").Value = 11592.286
$ws.Range("I52").Value = 11429.2
$ws.Range("K52").Value = 11429.2
$ws.Range("M52").Value = -11203.2
$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178
$ws.Range("H126").Value = 5473.231
$ws.Range("I126").Value = 5506
$ws.Range("J126").Value = 5435
$ws.Range("K126").Value = 16518
$ws.Range("L126").Value = 16305
$ws.Range("M126").Value = -14048
$ws.Range("N126").Value = -21245

